$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Update time_taken timestamps on the "data" sheet (F2:F6) ---
$ws.Cells.Item(2, 6).Value = "2021-10-05 14:35:32.421508"
$ws.Cells.Item(3, 6).Value = "2021-10-05 14:35:32.421516"
$ws.Cells.Item(4, 6).Value = "2021-10-05 14:35:32.421519"
$ws.Cells.Item(5, 6).Value = "2021-10-05 14:35:32.421522"
$ws.Cells.Item(6, 6).Value = "2021-10-05 14:35:32.421525"

# --- Add the new "metadata" sheet right after "data" ---
$newWs = $wb.Worksheets.Add($null, $ws)
$newWs.Name = "metadata"

# Copy header styling (bold + border + centered) from the "data" sheet's
# header row so the new header cells share the exact same style.
$ws.Range("B1").Copy()
$newWs.Range("B1:G1").PasteSpecial(-4122)

# Copy the style used for the numeric index column ("A") too.
$ws.Range("A2").Copy()
$newWs.Range("A2").PasteSpecial(-4122)

# --- Header row ---
$newWs.Cells.Item(1, 2).Value = "data_name"
$newWs.Cells.Item(1, 3).Value = "data_id"
$newWs.Cells.Item(1, 4).Value = "data_version"
$newWs.Cells.Item(1, 5).Value = "data_version_created"
$newWs.Cells.Item(1, 6).Value = "panel_query_time"
$newWs.Cells.Item(1, 7).Value = "panel_get_request"

# --- Data row ---
$newWs.Cells.Item(2, 1).Value = 0
$newWs.Cells.Item(2, 2).Value = "Renal Amyloidosis"
$newWs.Cells.Item(2, 3).Value = 191
# Force "0.21" to stay a text value (not coerced to the number 0.21), then
# reset the cell's style back to the default so no border/bold sneaks in.
$newWs.Cells.Item(2, 4).Value = "'0.21"
$newWs.Cells.Item(2, 4).Style = "Normal"
$newWs.Cells.Item(2, 5).Value = "2021-05-18T01:25:24.358074Z"
$newWs.Cells.Item(2, 6).Value = "2021-10-05 14:35:32.417669"
$newWs.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/191/?format=json"

$ws.Select()
